$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column AG: "Financial Secrecy Index (2018)" ---
# Header cell AG1 should look like the other header cells (e.g. AF1 / style 1).
$ws.Range("AF1").Copy()
$ws.Range("AG1").PasteSpecial(-4122)
$ws.Range("AG1").Value = "Financial Secrecy Index (2018)"

# Body cells AG2:AG28 get a red checkmark, matching the existing red-checkmark
# cells elsewhere in the sheet (e.g. G16 uses that style).
$ws.Range("G16").Copy()
$ws.Range("AG2:AG28").PasteSpecial(-4122)
$ws.Range("AG2:AG28").Value = "✓"

# --- Row 30 / cell A30: wrap text + taller row ---
$ws.Range("A30").WrapText = $true
$ws.Range("A30").RowHeight = 45.9

$excel.CutCopyMode = $false

# Leave the selection where the author left it when they saved the file.
$ws.Range("A28").Select() | Out-Null
